$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# R5 section gains two more method rows:
#  - "lookForVehicle" right after the existing
#    "lookForVehicles(displacement...)" row (old row 48), and the "Client"
#    row that used to be row 48 shifts down to row 49.
#  - "getInterestingVehicles" before the "searchVehicle(licensePlate...)"
#    row (old row 49, now at row 50 after the previous insert).
$ws.Rows.Item(48).Insert()
$ws.Rows.Item(50).Insert()

# Shared-string table order matters: the original author's edit appended
# "getInterestingVehicles" (row 50) before "lookForVehicle" (row 48), so we
# set the values in that same order to line up the resulting uniqueCount
# indices with the authored file.
$ws.Range("C50").Value = "getInterestingVehicles() : ArrayList<model.Vehicle>"
$ws.Range("C48").Value = "lookForVehicle(array : ArrayList<model.Vehicle>, vehicle : model.Vehicle) : boolean"

# New requirement block R6, appended after the (now shifted) Vehicle/toString
# row which sits at row 52.
$ws.Range("A53").Value = "R6. Show a complete report with all the data of the vehicles of interest of a client"
$ws.Range("B53").Value = "Main"
$ws.Range("C53").Value = "showVehiclesOfInterest() : void"

$ws.Range("C54").Value = "printVehicles() : int"

$ws.Range("B55").Value = "Business"
$ws.Range("C55").Value = "searchClient(id : int) : model.Client"

$ws.Range("B56").Value = "Client"
$ws.Range("C56").Value = "getInterestingVehicles() : ArrayList<model.Vehicle>"

# Keep the selection / top-left cell consistent with the new bottom of sheet.
$ws.Range("A57").Select() | Out-Null
